# Update carjacking-by-month-yoy report to include data through 2022-10-06
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the workbook title
$ws.Name = "Through 2022-10-06"

# Update the "October (through 10-05)" label to "October (through 10-06)"
$ws.Range("A11").Value = "October (through 10-06)"

# Update July 2022 value (I8)
$ws.Range("I8").Value = 163

# Update October row (row 11) with the latest data
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 9
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = 45
$ws.Range("I11").Value = 18

# Update Total row (row 12) to reflect new sums
$ws.Range("B12").Value = 232
$ws.Range("C12").Value = 438
$ws.Range("E12").Value = 563
$ws.Range("F12").Value = 429
$ws.Range("G12").Value = 934
$ws.Range("H12").Value = 1292
